# Updates the cryptocurrency price/volume table in-place to reflect the
# latest scrape (GitHub Actions job "Updated cryptos list").
# Column D = Price, Column E = Volume(1h). Row 41/42 (Bittensor/Filecoin)
# swapped rank position, so their Coin/Link/Price/Volume values are
# rewritten together.
#
# Numeric-looking Price strings are written with a leading apostrophe so
# Excel keeps them as literal text (matching the original inlineStr
# formatting) instead of auto-converting to numbers and dropping
# trailing zeros / thousands separators.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.094.27"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "2.531.00"
$ws.Range("E3").Value = "  +2.10%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'540.58"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").Value = "'143.66"
$ws.Range("E6").Value = "  -3.05%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "2.527.87"
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").Value = "'5.54"
$ws.Range("E12").Value = "  +3.92%  "
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").Value = "2.972.43"
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("D15").Value = "'23.56"
$ws.Range("E15").Value = "  -2.59%  "
$ws.Range("D16").Value = "59.074.31"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("D18").Value = "2.521.73"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "'11.19"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("E20").Value = "  -1.93%  "
$ws.Range("D21").Value = "'324.29"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("E22").Value = "  +3.10%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'61.80"
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("D25").Value = "'0.438"
$ws.Range("E25").Value = "  -5.22%  "
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("D27").Value = "'0.995"
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("D28").Value = "'7.92"
$ws.Range("E28").Value = "  +2.42%  "
$ws.Range("D29").Value = "0.0₃0776"
$ws.Range("E30").Value = "  -1.57%  "
$ws.Range("D31").Value = "'6.64"
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("E32").Value = "  -8.95%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").Value = "'1.46"
$ws.Range("E34").Value = "  +6.22%  "
$ws.Range("D35").Value = "'158.20"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").Value = "'18.64"
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("E37").Value = "  -4.23%  "
$ws.Range("D38").Value = "'1.61"
$ws.Range("E38").Value = "  -7.43%  "
$ws.Range("D39").Value = "'36.98"
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("E40").Value = "  -5.03%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'3.70"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "'295.04"
$ws.Range("E42").Value = "  -7.76%  "
$ws.Range("D43").Value = "'0.821"
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("D44").Value = "'0.997"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("E45").Value = "  +2.53%  "
$ws.Range("D46").Value = "'10.79"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("E47").Value = "  -1.66%  "
$ws.Range("D48").Value = "'18.58"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").Value = "'122.33"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("E50").Value = "  -2.26%  "
$ws.Range("E51").Value = "  -0.98%  "

